$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19, shifting existing rows 19..135 down to 20..136.
$ws.Rows(19).Insert()

# Populate the newly inserted row 19 with its data (same record shape as the
# row that used to be there, but with its own date / price / origin values).
$ws.Range("A19").Value2 = 6
$ws.Range("B19").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C19").Value2 = "Metropolitana"
$ws.Range("D19").Value2 = 44561
$ws.Range("E19").Value2 = 13
$ws.Range("F19").Value2 = "Fruta"
$ws.Range("G19").Value2 = 100101
$ws.Range("H19").Value2 = "Berries"
$ws.Range("I19").Value2 = 100101004
$ws.Range("J19").Value2 = "Frambuesa"
$ws.Range("K19").Value2 = "Sin especificar"
$ws.Range("L19").Value2 = "Primera"
$ws.Range("M19").Value2 = 250
$ws.Range("N19").Value2 = 6000
$ws.Range("O19").Value2 = 7000
$ws.Range("P19").Value2 = 6500
$ws.Range("Q19").Value2 = "`$/bandeja 2 kilos"
$ws.Range("R19").Value2 = "Provincia de Colchagua"
$ws.Range("S19").Value2 = 3250
$ws.Range("T19").Value2 = 2
